$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2052.923
$ws.Range("J17").Value = 2149
$ws.Range("L17").Value = 6447
$ws.Range("N17").Value = -6783

$ws.Range("H19").Value = 1045.3334
$ws.Range("I19").Value = 931.875
$ws.Range("K19").Value = 931.875
$ws.Range("M19").Value = -756.875

$ws.Range("H103").Value = 962.6111
$ws.Range("I103").Value = 596
$ws.Range("K103").Value = 1788
$ws.Range("M103").Value = -1202

$ws.Range("H107").Value = 521.25
$ws.Range("I107").Value = 600.55554
$ws.Range("K107").Value = 600.55554
$ws.Range("M107").Value = 1319.44446

$ws.Range("H112").Value = 2223
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2223
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 6669
$ws.Range("N112").Value = -8885
$ws.Range("M112").ClearContents()

$ws.Range("H115").Value = 760.94446
$ws.Range("I115").Value = 760.94446
$ws.Range("K115").Value = 2282.83338
$ws.Range("M115").Value = -715.83338

$ws.Range("H132").Value = 4269.2705
$ws.Range("I132").Value = 4468.147
$ws.Range("K132").Value = 13404.441
$ws.Range("M132").Value = -10874.441

$ws.Range("H134").Value = 31838.928
$ws.Range("J134").Value = 31838.928
$ws.Range("L134").Value = 31838.928
$ws.Range("N134").Value = -41978.928

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6034.1045
$ws.Range("I32").Value = 4875.656
$ws.Range("K32").Value = 4875.656
$ws.Range("M32").Value = -4588.656

$ws.Range("H61").Value = 2577.3044
$ws.Range("I61").Value = 820.65
$ws.Range("J61").Value = 14288.333
$ws.Range("K61").Value = 820.65
$ws.Range("L61").Value = 14288.333
$ws.Range("M61").Value = -608.65
$ws.Range("N61").Value = -14712.333

$ws.Range("H97").Value = 2297.5557
$ws.Range("I97").Value = 446.38235
$ws.Range("K97").Value = 446.38235
$ws.Range("M97").Value = 49.61765000000003

$ws.Range("H110").Value = 329.17648
$ws.Range("J110").Value = 322.4
$ws.Range("L110").Value = 322.4
$ws.Range("N110").Value = -4412.4

$ws.Range("H136").Value = 2577.3044
$ws.Range("I136").Value = 820.65
$ws.Range("J136").Value = 14288.333
$ws.Range("K136").Value = 2461.95
$ws.Range("L136").Value = 42864.999
$ws.Range("M136").Value = 88.05000000000018
$ws.Range("N136").Value = -47964.999

$ws.Range("H138").Value = 52126.855
$ws.Range("J138").Value = 52126.855
$ws.Range("L138").Value = 52126.855
$ws.Range("N138").Value = -62406.855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 9447.091
$ws.Range("I94").Value = 452.375
$ws.Range("K94").Value = 452.375
$ws.Range("M94").Value = -1.375

$ws.Range("H134").Value = 7247.961
$ws.Range("I134").Value = 7256.4595
$ws.Range("K134").Value = 21769.3785
$ws.Range("M134").Value = -19234.3785

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1875.0769
$ws.Range("I16").Value = 1307.1818
$ws.Range("K16").Value = 1307.1818
$ws.Range("M16").Value = -1020.1818

$ws.Range("H107").Value = 553.1429000000001
$ws.Range("I107").Value = 311.8
$ws.Range("K107").Value = 311.8
$ws.Range("M107").Value = 1608.2

$ws.Range("H113").Value = 1875.0769
$ws.Range("I113").Value = 1307.1818
$ws.Range("K113").Value = 1307.1818
$ws.Range("M113").Value = 862.8181999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 482.8
$ws.Range("J5").Value = 537.5
$ws.Range("L5").Value = 1612.5
$ws.Range("N5").Value = -1836.5

$ws.Range("H23").Value = 216.71428
$ws.Range("I23").Value = 141.66667
$ws.Range("J23").Value = 273
$ws.Range("K23").Value = 425.00001
$ws.Range("L23").Value = 819
$ws.Range("M23").Value = -190.00001
$ws.Range("N23").Value = -1289

$ws.Range("H97").Value = 387.6
$ws.Range("J97").Value = 387
$ws.Range("L97").Value = 1161
$ws.Range("N97").Value = -2153

$ws.Range("H101").Value = 25401.934
$ws.Range("J101").Value = 28002.637
$ws.Range("L101").Value = 84007.91099999999
$ws.Range("N101").Value = -88875.91099999999

$ws.Range("H132").Value = 4163.8945
$ws.Range("I132").Value = 3336.75
$ws.Range("K132").Value = 30030.75
$ws.Range("M132").Value = -27500.75

$ws.Range("H135").Value = 482.8
$ws.Range("J135").Value = 537.5
$ws.Range("L135").Value = 4837.5
$ws.Range("N135").Value = -9907.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6405.6665
$ws.Range("I126").Value = 5606.4287
$ws.Range("J126").Value = 7524.6
$ws.Range("K126").Value = 16819.2861
$ws.Range("L126").Value = 22573.8
$ws.Range("M126").Value = -14349.2861
$ws.Range("N126").Value = -27513.8

$ws.Range("H136").Value = 26467.88
$ws.Range("J136").Value = 26467.88
$ws.Range("L136").Value = 79403.64
$ws.Range("N136").Value = -84503.64

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 667.2857
$ws.Range("J16").Value = 810.55554
$ws.Range("L16").Value = 810.55554
$ws.Range("N16").Value = -1150.55554

$ws.Range("H61").Value = 1789.3043
$ws.Range("I61").Value = 1707.75
$ws.Range("J61").Value = 2333
$ws.Range("K61").Value = 1707.75
$ws.Range("L61").Value = 2333
$ws.Range("M61").Value = -1505.75
$ws.Range("N61").Value = -2737

$ws.Range("H113").Value = 1789.3043
$ws.Range("I113").Value = 1707.75
$ws.Range("J113").Value = 2333
$ws.Range("K113").Value = 1707.75
$ws.Range("L113").Value = 2333
$ws.Range("M113").Value = 462.25
$ws.Range("N113").Value = -6673

$ws.Range("H132").Value = 4022.2727
$ws.Range("J132").Value = 4164.625
$ws.Range("L132").Value = 12493.875
$ws.Range("N132").Value = -17553.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 11676.2
$ws.Range("I23").Value = 999
$ws.Range("J23").Value = 14345.5
$ws.Range("K23").Value = 999
$ws.Range("L23").Value = 14345.5
$ws.Range("M23").Value = -770
$ws.Range("N23").Value = -14803.5

$ws.Range("H107").Value = 5031.174
$ws.Range("I107").Value = 7535.4287
$ws.Range("K107").Value = 22606.2861
$ws.Range("M107").Value = -20686.2861

$ws.Range("H122").Value = 102498.75
$ws.Range("I122").Value = 103332.336
$ws.Range("K122").Value = 309997.008
$ws.Range("M122").Value = -307547.008

$ws.Range("H126").Value = 5051.3335
$ws.Range("I126").Value = 4355
$ws.Range("J126").Value = 5399.5
$ws.Range("K126").Value = 13065
$ws.Range("L126").Value = 16198.5
$ws.Range("M126").Value = -10595
$ws.Range("N126").Value = -21138.5

$ws.Range("H132").Value = 2676.55
$ws.Range("I132").Value = 2467.7646
$ws.Range("K132").Value = 7403.293799999999
$ws.Range("M132").Value = -4873.293799999999

$ws.Range("H136").Value = 5878.2666
$ws.Range("J136").Value = 7802
$ws.Range("L136").Value = 23406
$ws.Range("N136").Value = -28506
